# Fix date errors in the original Dec 2019 finance data.
# Several rows were mistakenly entered with 2020 serial dates instead of
# the correct 2019 dates (exactly one year off -> 366 days, since 2020 was
# a leap year). Correct them in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3: 1-Dec-2020 (44166) -> 1-Dec-2019 (43800)
$ws.Range("A3").Value = 43800

# A23: 13-Dec-2020 (44178) -> 13-Dec-2019 (43812)
$ws.Range("A23").Value = 43812

# A52: 28-Dec-2020 (44193) -> 28-Dec-2019 (43827)
$ws.Range("A52").Value = 43827

# A53: 26-Dec-2020 (44191) -> 26-Dec-2019 (43825)
$ws.Range("A53").Value = 43825

# Update the view: scroll the window down and move the active selection
# to A56, matching where the user was working after the fix.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A56").Select()
